# "Generate Report for Handoff"
# Updates the localization-status report: the Status column moves from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# associated handoff timestamps are refreshed. The Status columns are
# narrowed to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff timestamps ---
$overview.Range("G2").Value = "2016-08-18 04:59:25"
$dede.Range("H2").Value = "2016-08-18 04:59:25"
$zhcn.Range("H2").Value = "2016-08-18 04:59:20"

# --- Narrow the now-shorter Status columns to fit the new text ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
